# Sprint 4 Development and Issues fixed
# Localize the Sheet1 header cells from English to Spanish:
#   A1 (shared string "Key Word") -> "Palabra clave"
#   B1 (shared string "Quantity") -> "Cantidad"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Palabra clave"
$ws.Range("B1").Value = "Cantidad"
